$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relocate the existing "IDB" source block from rows 24-25 down to rows 30-31 ---
# (the new MSME size-class table is being inserted where rows 24-25 used to be)
$ws.Range("A30").Value = $ws.Range("A24").Value2
$ws.Range("A30").Style = "title"
$ws.Range("A31").Value = $ws.Range("A25").Value2
$ws.Range("A31").Style = "source"

# --- New table header (row 21): Number of employees / Assets / Turnover ---
$ws.Range("B21").Value = "Number of employees"
$ws.Range("B21").Style = "title"
$ws.Range("C21").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C21").Style = "title"
$ws.Range("D21").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D21").Style = "title"

# --- Data rows 22-25: Micro / Small / Medium / Large, values left blank ---
$ws.Range("A22").Value = "Micro"
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""

$ws.Range("A23").Value = "Small"
$ws.Range("B23").Value = ""
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = ""

$ws.Range("A24").Value = "Medium"
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""

$ws.Range("A25").Value = "Large"
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""
